# Fix systematic spacing issue between header bar and body text
# (content changes: consolidate CORE COMPETENCIES bullets into one summary
#  line, and add a detailed TECHNICAL SKILLS section near the end.)

$d = $word.ActiveDocument
$bullet = [char]0x2022

# ---------------------------------------------------------------------
# 1) CORE COMPETENCIES: collapse the three detailed bullet paragraphs
#    into a single short summary paragraph.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Survey Methodology & Research Design: Survey Design and Questionnaire Development for Political and Market Research " + $bullet + " Sampling Methodology and Statistical Analysis (R, SPSS, Stata, OSCAR) " + $bullet + " Random Device Engagement (RDE), Text Message, Web Panel, and Live Telephone Calling " + $bullet + " Expert Testimony and Consultation on Research Methodology",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Survey Methodology & Research Design " + $bullet + " Redistricting & Geospatial Analysis " + $bullet + " Data Analysis & Visualization",
    2)

$d.Content.Find.Execute(
    "Redistricting & Geospatial Analysis: Redistricting Software Development and Boundary Estimation Systems " + $bullet + " Geospatial Analysis: ArcGIS, Quantum GIS, GRASS, OSGeo, PostGIS " + $bullet + " Choropleths and Hexagonal Grid Maps for Demographic Visualization " + $bullet + " Court Case Analysis and Expert Testimony for Redistricting",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "",
    2)

$d.Content.Find.Execute(
    "Data Analysis & Visualization: Advanced Statistical Modeling and Analysis (Regression, Clustering, Segmentation) " + $bullet + " Data Visualization: Tableau, PowerBI, Seaborn, Matplotlib, d3.js " + $bullet + " Consumer Behavior Analysis and Market Segmentation " + $bullet + " Multi-million Dollar Research Project Management",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "",
    2)

# Remove the two now-empty paragraphs left behind by the replacements above
# (an "empty" paragraph's Range.Text is just the lone paragraph-mark, "\r").
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim() -eq "") {
        $p.Range.Delete()
        $i = $i - 1
    }
}

# ---------------------------------------------------------------------
# 2) Append a new "TECHNICAL SKILLS" section (Heading2 + 3 paragraphs)
#    right before the closing "For a more detailed..." paragraph.
# ---------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "Developed advanced segmentation models using demographic, psychographic, and behavioral data") {
        $target = $p
    }
}

$target.Range.InsertParagraphAfter()
$idx = $target.Index + 1
$d.Paragraphs.Item($idx).Range.Text = "TECHNICAL SKILLS"

$d.Paragraphs.Item($idx).Range.InsertParagraphAfter()
$idx2 = $idx + 1
$d.Paragraphs.Item($idx2).Range.Text = "SURVEY METHODOLOGY & RESEARCH DESIGN Survey Design and Questionnaire Development for Political and Market Research; Sampling Methodology and Statistical Analysis (R, SPSS, Stata, OSCAR); Random Device Engagement (RDE), Text Message, Web Panel, and Live Telephone Calling; Expert Testimony and Consultation on Research Methodology"

$d.Paragraphs.Item($idx2).Range.InsertParagraphAfter()
$idx3 = $idx2 + 1
$d.Paragraphs.Item($idx3).Range.Text = "REDISTRICTING & GEOSPATIAL ANALYSIS Redistricting Software Development and Boundary Estimation Systems; Geospatial Analysis; Choropleths and Hexagonal Grid Maps for Demographic Visualization; Court Case Analysis and Expert Testimony for Redistricting"

$d.Paragraphs.Item($idx3).Range.InsertParagraphAfter()
$idx4 = $idx3 + 1
$d.Paragraphs.Item($idx4).Range.Text = "DATA ANALYSIS & VISUALIZATION Advanced Statistical Modeling and Analysis (Regression, Clustering, Segmentation); Data Visualization; Consumer Behavior Analysis and Market Segmentation; Multi-million Dollar Research Project Management"

# Apply the Heading2 style to the new section title only AFTER all the
# InsertParagraphAfter calls are done (otherwise the heading style leaks
# into the freshly-inserted sibling paragraphs).
$d.Paragraphs.Item($idx).Style = "Heading2"
